$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0.00275442249619973
$ws.Range("E2").Value = 0.9514223331899243
$ws.Range("F2").Value = 0.003188848495483398
$ws.Range("G2").Value = 0.9960679680569756

# Row 3
$ws.Range("B3").Value = 11
$ws.Range("C3").Value = 0.05158144172721063
$ws.Range("E3").Value = 0.09029711554240791
$ws.Range("F3").Value = 0.01740074157714844
$ws.Range("G3").Value = 0.9758368432954937

# Row 4
$ws.Range("B4").Value = 16
$ws.Range("C4").Value = 0.0566794055335246
$ws.Range("E4").Value = 0.0003881827136121285
$ws.Range("F4").Value = 0.02093791961669922
$ws.Range("G4").Value = 0.9669067670575243

# Row 5
$ws.Range("B5").Value = 4371
$ws.Range("C5").Value = 0.05675808555049239
$ws.Range("E5").Value = 0.0009994372517873483
$ws.Range("F5").Value = 2.680588722229004
$ws.Range("G5").Value = 0.6006062222089232

# Row 6
$ws.Range("B6").Value = 21724
$ws.Range("C6").Value = 0.05670708231967949
$ws.Range("E6").Value = [double]"9.993183601156125E-05"
$ws.Range("F6").Value = 13.74774575233459
$ws.Range("G6").Value = 0.3598264157717663

# Row 7
$ws.Range("B7").Value = 82212
$ws.Range("C7").Value = 0.05670198197752034
$ws.Range("E7").Value = [double]"9.98095797220148E-06"
$ws.Range("F7").Value = 53.45275068283081
$ws.Range("G7").Value = 0.1776273892192539

# Row 8
$ws.Range("B8").Value = 190229
$ws.Range("C8").Value = 0.05670147271032179
$ws.Range("E8").Value = [double]"9.993974735381498E-07"
$ws.Range("F8").Value = 118.4360136985779
$ws.Range("G8").Value = 0.09369313456201694
